$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") stores the "last changed" date for every data row.
# Find the last used row (based on column A, "Beteckning") so all data rows are covered.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 339 }

$ws.Range("C2:C" + $lastRow).Value = 45175
